# Generate Report for Handoff
# Rename the localized markdown/xliff identifiers from
#   6a7d3b36-09c3-4ea7-9c27-0941ee66d382  ->  ad524f1c-9d01-4152-ac78-f7db533094fb
# and bump the handoff/generate timestamps, across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "6a7d3b36-09c3-4ea7-9c27-0941ee66d382"
$newGuid = "ad524f1c-9d01-4152-ac78-f7db533094fb"

$newZhXlf = "$newGuid.703042b5688bf47aa2cf97d1f0a04d1ec3831e5d.zh-cn.xlf"
$newDeXlf = "$newGuid.703042b5688bf47aa2cf97d1f0a04d1ec3831e5d.de-de.xlf"

$hyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d09cde0e2ca726a16af6cc9af4c88ddaa4961aab/e2e/$oldGuid.md"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Sheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-26 06:56:15"

# Recreate the hyperlink on B2 so its display text matches the new cell text
# (setting a property on an existing Hyperlink object duplicates the entry in
# this runtime, so clear the collection first and re-add it cleanly).
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkUrl, [Type]::Missing, [Type]::Missing, "e2e\$newGuid.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Sheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = "2016-08-26 06:56:11"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkUrl, [Type]::Missing, [Type]::Missing, "$newGuid.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Sheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = $newDeXlf
$wsDeDe.Range("H2").Value = "2016-08-26 06:56:15"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkUrl, [Type]::Missing, [Type]::Missing, "$newGuid.md") | Out-Null
